$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.74"
$ws.Range("E2").Value = "'1.68%"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'2.93%"
$ws.Range("D4").Value = "'5.605"
$ws.Range("E4").Value = "'-4.74%"
$ws.Range("D5").Value = "'0.08165"
$ws.Range("E5").Value = "'1.68%"
$ws.Range("D6").Value = "'2.038"
$ws.Range("E6").Value = "'5.79%"
$ws.Range("D7").Value = "'8.743"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("D8").Value = "'4.536"
$ws.Range("E8").Value = "'-1.18%"
$ws.Range("D9").Value = "'3.000"
$ws.Range("E9").Value = "'1.92%"
$ws.Range("D10").Value = "'0.9189"
$ws.Range("E10").Value = "'-1.55%"
$ws.Range("D11").Value = "'0.1260"
$ws.Range("E11").Value = "'-1.02%"
$ws.Range("E12").Value = "'-0.62%"
$ws.Range("D13").Value = "'0.09314"
$ws.Range("E13").Value = "'2.09%"
$ws.Range("D14").Value = "'0.03745"
$ws.Range("E14").Value = "'5.64%"
$ws.Range("D15").Value = "'0.1055"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("D16").Value = "'0.001301"
$ws.Range("E16").Value = "'0.34%"
$ws.Range("D17").Value = "'0.006305"
$ws.Range("E17").Value = "'1.29%"
$ws.Range("E19").Value = "'-2.23%"
$ws.Range("D20").Value = "'8.471"
$ws.Range("E20").Value = "'-2.93%"
$ws.Range("E21").Value = "'-1.84%"
$ws.Range("D22").Value = "'0.2517"
$ws.Range("E22").Value = "'4.46%"
$ws.Range("D23").Value = "'0.04419"
$ws.Range("E23").Value = "'0.17%"
$ws.Range("D24").Value = "'0.001265"
$ws.Range("E24").Value = "'0.23%"
$ws.Range("D25").Value = "'0.004307"
$ws.Range("E25").Value = "'-2.20%"
$ws.Range("D26").Value = "'0.0001183"
$ws.Range("E26").Value = "'3.80%"
$ws.Range("D39").Value = "'0.02741"
$ws.Range("E39").Value = "'12.37%"
$ws.Range("E40").Value = "'3.00%"
$ws.Range("D41").Value = "'0.007664"
$ws.Range("E41").Value = "'3.47%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1410"
$ws.Range("E42").Value = "'0.35%"
$ws.Range("B43").Value = "Dexo"
$ws.Range("C43").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D43").Value = "'0.009154"
$ws.Range("E43").Value = "'-3.56%"
$ws.Range("D44").Value = "'0.002249"
$ws.Range("E44").Value = "'6.12%"
$ws.Range("D45").Value = "'0.01127"
$ws.Range("E45").Value = "'13.25%"
$ws.Range("D46").Value = "'0.00006893"
$ws.Range("E46").Value = "'2.28%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.26%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002284"
$ws.Range("E48").Value = "'60.56%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003537"
$ws.Range("E49").Value = "'17.89%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.26%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.26%"
